# Actualización automática hashcode lun sep 21 01:27:27 CEST 2020
#
# This script updates the "hashcode" column (column B) for a set of rows
# identified by their key in column A, replacing the old hash value with
# a freshly computed one, mirroring the upstream automated commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of key (column A value) -> new hash (column B value) to apply.
$updates = @(
    @{ Key = "05-050316TC";       New = "2a8524da19a261ecdef6891100f68859" },
    @{ Key = "05-050104A";        New = "f7945b435d376f43969ae850a7cc68cb" },
    @{ Key = "05-050101A";        New = "45345d734b099da46e786c83e8f28c96" },
    @{ Key = "05-050102A";        New = "74c498ae62afc36eaf69fb2be262b624" },
    @{ Key = "05-050203TC";       New = "d321d6ac78ad3d5397984693326c7178" },
    @{ Key = "05-050009TC";       New = "b2c2d7b0c6e1e482e2baebfaa3e80238" },
    @{ Key = "05-050205TP";       New = "811e4b110a2cffba77fce045c7017d73" },
    @{ Key = "05-050009TP";       New = "67e8de9238b1d980854c534789e8446c" },
    @{ Key = "05-050205TC";       New = "869c621bbced2dd1e9009bcaac137d49" },
    @{ Key = "05-050101TP";       New = "beba7bce29c4068483cd10898052ff4a" },
    @{ Key = "05-050201TC";       New = "f23660b688dfd8a0463a2ff716f4e132" },
    @{ Key = "05-050005TP";       New = "7d3192fea74a6be1ead9e53c83c35f0f" },
    @{ Key = "05-0709-070905BTC"; New = "0841f66eec1f7caf51680bed6f5054c6" },
    @{ Key = "05-050205A";        New = "7c7e26fef28b133513b0e1d817db11ed" },
    @{ Key = "05-050208TP";       New = "3bb24bf20af84bd73d4fd48e30da03f3" },
    @{ Key = "05-050009A";        New = "46abcc7d85f2732d753478da077c6dad" },
    @{ Key = "05-050201A";        New = "94c8a699ba72fa2ba49483e62eaeeb5b" },
    @{ Key = "05-050004A";        New = "309f583d917950c45f020d6995e0ecb3" },
    @{ Key = "05-050308TP";       New = "2dfdedb2c6659147cc3aefedac967c38" },
    @{ Key = "05-050005A";        New = "0500c3294f2fe90971052abfee60871b" },
    @{ Key = "05-050208A";        New = "2ede366eee4394e48ea0925f9464345c" },
    @{ Key = "05-050206TP";       New = "87f7d8c8d5f14748512c9245c79f6ea6" },
    @{ Key = "05-050206TC";       New = "e992428de39ad6cc52cb72f089587295" },
    @{ Key = "05-050304TC";       New = "c73244e4d02da93b2f5418460dd36c9d" },
    @{ Key = "05-050206A";        New = "d174fa8fbca0c777f41402c2571309ad" },
    @{ Key = "05-050315A";        New = "14cb8d34718c47516b19ad2970bcf17c" },
    @{ Key = "05-050315TP";       New = "73dcb4033cf74069e3da205ee99500a5" },
    @{ Key = "05-050102TP";       New = "0a647b4a3f32e50bca26867df944df5e" },
    @{ Key = "05-050104TC";       New = "369163dccc3c430a954a07963037cfd1" },
    @{ Key = "05-050104TM";       New = "55ee70e9919cf8142a528225a340560d" },
    @{ Key = "05-050104TP";       New = "e8dfad8ff97156163b1440cb8b6475c6" }
)

# Determine the used range of column A so we know how many rows to scan.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Build a lookup of key -> row number by scanning column A once.
# NOTE: use .Value2 (not the parameterless .Value get) for reads here --
# the host's plain `.Value` property read returns the member signature
# string instead of the cell's contents; `.Value2`/`.Text` read correctly,
# and plain `.Value = ...` assignment (write) is unaffected.
$keyToRow = @{}
for ($r = 1; $r -le $lastRow; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    if ($null -ne $a) {
        $keyToRow[[string]$a] = $r
    }
}

foreach ($update in $updates) {
    $row = $keyToRow[$update.Key]
    if ($null -ne $row) {
        $ws.Cells.Item($row, 2).Value = $update.New
    }
}
